# Generate Report for Handback
# Refresh the handoff/handback timestamps for the 067e103b...md file
# (both the zh-cn and de-de detail sheets) and roll the refreshed
# de-de handoff timestamp up into the Overview sheet's
# "Latest HO Xliff Generate Date" column for that row.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")
$overview = $wb.Worksheets.Item("Overview")

# zh-cn: row 2 corresponds to 067e103b-4096-4abf-beff-6194a708d5d1.md
$zhcn.Range("H2").Value = "2016-08-28 00:45:53"
$zhcn.Range("K2").Value = "2016-08-28 00:46:13"

# de-de: row 2 corresponds to 067e103b-4096-4abf-beff-6194a708d5d1.md
$dede.Range("H2").Value = "2016-08-28 00:45:57"
$dede.Range("K2").Value = "2016-08-28 00:46:19"

# Overview: G column ("Latest HO Xliff Generate Date") mirrors the
# de-de handoff datetime for each file.
$overview.Range("G2").Value = "2016-08-28 00:45:57"
$overview.Range("G3").Value = "2016-08-28 00:44:52"
